# Updates cryptos list cell values per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.140.27"
$ws.Range("E2").Value = "  +2.92%  "

$ws.Range("D3").Value = "2.058.09"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'230.02"
$ws.Range("E5").Value = "  +1.88%  "

$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  +2.24%  "

$ws.Range("D7").Value = "'59.54"
$ws.Range("E7").Value = "  +7.90%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +3.35%  "

$ws.Range("D10").Value = "'0.0812"
$ws.Range("E10").Value = "  +4.66%  "

$ws.Range("E11").Value = "  +2.31%  "

$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'14.74"
$ws.Range("E12").Value = "  +5.44%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.362.28"
$ws.Range("E13").Value = "  +2.37%  "

$ws.Range("D14").Value = "'21.21"
$ws.Range("E14").Value = "  +7.65%  "

$ws.Range("D15").Value = "'0.757"
$ws.Range("E15").Value = "  +3.11%  "

$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").Value = "2.054.25"
$ws.Range("E17").Value = "  +2.17%  "

$ws.Range("D18").Value = "38.005.83"
$ws.Range("E18").Value = "  +2.80%  "

$ws.Range("E19").Value = "  +1.78%  "

$ws.Range("D20").Value = "'69.88"
$ws.Range("E20").Value = "  +2.49%  "

$ws.Range("E21").Value = "  +3.23%  "

$ws.Range("D22").Value = "'225.15"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").Value = "'0.995"
$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("E25").Value = "  +4.36%  "

$ws.Range("D26").Value = "'166.45"
$ws.Range("E26").Value = "  +1.47%  "

$ws.Range("D27").Value = "'9.25"
$ws.Range("E27").Value = "  +4.08%  "

$ws.Range("E28").Value = "  +7.45%  "

$ws.Range("D29").Value = "'19.03"

$ws.Range("D30").Value = "'1.33"
$ws.Range("E30").Value = "  +2.58%  "

$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  +3.02%  "

$ws.Range("E32").Value = "  +3.95%  "

$ws.Range("E33").Value = "  +3.23%  "

$ws.Range("E34").Value = "  +10.61%  "

$ws.Range("D35").Value = "'0.0607"
$ws.Range("E35").Value = "  +1.38%  "

$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("D37").Value = "'6.18"
$ws.Range("E37").Value = "  +15.75%  "

$ws.Range("E38").Value = "  +5.71%  "

$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "1.530.21"
$ws.Range("E40").Value = "  +4.87%  "

$ws.Range("D41").Value = "'98.24"
$ws.Range("E41").Value = "  +3.91%  "

$ws.Range("D42").Value = "'17.12"
$ws.Range("E42").Value = "  +7.48%  "

$ws.Range("D43").Value = "'0.0216"
$ws.Range("E43").Value = "  +2.60%  "

$ws.Range("D44").Value = "'2.87"
$ws.Range("E44").Value = "  +4.50%  "

$ws.Range("D45").Value = "'0.0923"
$ws.Range("E45").Value = "  +1.96%  "

$ws.Range("E46").Value = "  +1.94%  "

$ws.Range("D47").Value = "'4.11"
$ws.Range("E47").Value = "  -3.58%  "

$ws.Range("E48").Value = "  +2.66%  "

$ws.Range("E49").Value = "  +3.00%  "

$ws.Range("D50").Value = "'7.12"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").Value = "2.250.30"
$ws.Range("E51").Value = "  +2.57%  "
